$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.125.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.863.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4677"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2829"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06574"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.881.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.062"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6694"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.141.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.127.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.367"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007242"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.00%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.150"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.321"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.950"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.373"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09655"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.376"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.466"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.094"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04656"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7004"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.086"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.718"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01856"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.406"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.511"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8596"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.941"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.002"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4163"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "984.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.187"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.078"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.73%  "
